$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from row 19 down to row 20 (so date style etc. carries over)
$ws.Range("A19:E19").Copy() | Out-Null
$ws.Range("A20:E20").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("A20").Value = 45986
$ws.Range("B20").Value = 2025
$ws.Range("C20").Value = 2.46481303148316
$ws.Range("D20").Value = 2026
$ws.Range("E20").Value = 3.633434696013671
